$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.995.25"
$ws.Range("E2").Value = "  +0.57%  "
$ws.Range("D3").Value = "2.291.38"
$ws.Range("E3").Value = "  +1.90%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "'112.32"
$ws.Range("E5").Value = "  -3.41%  "
$ws.Range("D6").Value = "'309.51"
$ws.Range("E6").Value = "  +2.16%  "
$ws.Range("D7").Value = "'0.633"
$ws.Range("E7").Value = "  +0.78%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").Value = "'0.616"
$ws.Range("E9").Value = "  -1.51%  "
$ws.Range("D10").Value = "'44.39"
$ws.Range("E10").Value = "  -4.61%  "
$ws.Range("D11").Value = "'0.0929"
$ws.Range("E11").Value = "  -0.81%  "
$ws.Range("D12").Value = "'55.26"
$ws.Range("E12").Value = "  -2.96%  "
$ws.Range("D13").Value = "'8.83"
$ws.Range("E13").Value = "  -3.77%  "
$ws.Range("E14").Value = "  +20.57%  "
$ws.Range("E15").Value = "  -0.70%  "
$ws.Range("D16").Value = "'15.54"
$ws.Range("E16").Value = "  +0.63%  "
$ws.Range("D17").Value = "2.637.88"
$ws.Range("E17").Value = "  +1.95%  "
$ws.Range("D18").Value = "2.287.86"
$ws.Range("E18").Value = "  +0.93%  "
$ws.Range("D19").Value = "42.911.22"
$ws.Range("E19").Value = "  +0.09%  "
$ws.Range("E20").Value = "  -0.53%  "
$ws.Range("D21").Value = "'7.21"
$ws.Range("E21").Value = "  -6.91%  "
$ws.Range("D22").Value = "'76.29"
$ws.Range("E22").Value = "  +2.87%  "
$ws.Range("D23").Value = "'3.53"
$ws.Range("E23").Value = "  -3.72%  "
$ws.Range("D24").Value = "'2.46"
$ws.Range("E24").Value = "  +3.91%  "
$ws.Range("D25").Value = "'256.68"
$ws.Range("E25").Value = "  +9.91%  "
$ws.Range("E26").Value = "  -4.56%  "
$ws.Range("E27").Value = "  -4.16%  "
$ws.Range("D28").Value = "'0.999"
$ws.Range("E28").Value = "  -0.15%  "
$ws.Range("D29").Value = "'39.06"
$ws.Range("E29").Value = "  -3.52%  "
$ws.Range("E30").Value = "  -0.14%  "
$ws.Range("D31").Value = "'22.38"
$ws.Range("E31").Value = "  +4.59%  "
$ws.Range("D32").Value = "'173.50"
$ws.Range("E32").Value = "  -1.33%  "
$ws.Range("E33").Value = "  -2.83%  "
$ws.Range("D34").Value = "'0.0902"
$ws.Range("E34").Value = "  -1.28%  "
$ws.Range("D35").Value = "'5.75"
$ws.Range("E35").Value = "  +1.13%  "
$ws.Range("E36").Value = "  +6.31%  "
$ws.Range("E37").Value = "  +1.10%  "
$ws.Range("D38").Value = "'4.17"
$ws.Range("E38").Value = "  -9.61%  "
$ws.Range("D39").Value = "'0.0377"
$ws.Range("E39").Value = "  +0.11%  "
$ws.Range("E40").Value = "  -1.80%  "
$ws.Range("D41").Value = "'2.60"
$ws.Range("E41").Value = "  -0.30%  "
$ws.Range("D42").Value = "'72.13"
$ws.Range("E42").Value = "  -1.04%  "
$ws.Range("E43").Value = "  -3.77%  "
$ws.Range("D44").Value = "'1.44"
$ws.Range("E44").Value = "  +6.84%  "
$ws.Range("E45").Value = "  -0.20%  "
$ws.Range("D46").Value = "'12.40"
$ws.Range("E46").Value = "  -8.52%  "
$ws.Range("E47").Value = "  +2.04%  "
$ws.Range("D48").Value = "'108.68"
$ws.Range("E48").Value = "  +0.70%  "
$ws.Range("D49").Value = "'8.92"
$ws.Range("E49").Value = "  +2.97%  "
$ws.Range("E50").Value = "  -4.72%  "
$ws.Range("D51").Value = "'0.0986"
$ws.Range("E51").Value = "  -0.55%  "
